# Error Calculations and Plots
# Apply data corrections to the missing_data worksheet:
#  - Fix a handful of individual cell values/blanks in rows 5-25
#  - Rows 26-35 (RM 232 .. SC 232) had data shifted: the RM 232 row is
#    removed, SC 5..SC 232 rows are corrected, and the final two rows
#    (old SC 193 / SC 232 duplicated rows) are deleted so the table
#    shrinks from A1:F35 to A1:F33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Isolated single-cell fixes (rows 5-25) ---
$ws.Range("D5").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("D11").Value = -15.5
$ws.Range("C19").Value = 13.2
$ws.Range("D19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9
$ws.Range("F24").Value = 16.78
$ws.Range("D25").Value = -15.5

# --- Rows 26-33 get fully corrected content ---
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 11.2
$ws.Range("D29").ClearContents()
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").ClearContents()

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").ClearContents()

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

# --- Remove the now-duplicated trailing rows so the table ends at row 33 ---
$ws.Rows("34:35").Delete()
